# "Les nouveaux slides, les nouvelles positions"
# 1) Move slide 2 ("Groupe 3" decorative slide) to position 4.
# 2) Refresh the cached date-field text (7/1/2021 -> 11/1/2021) on the
#    slide master and every slide layout's Date placeholder.

$p = $ppt.ActivePresentation

# --- 1) Reorder slides: slide at position 2 moves to position 4 ---
$s = $p.Slides.Item(2)
$s.MoveTo(4)

# --- 2) Update the cached date placeholder text everywhere it appears ---

# Slide master
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "11/1/2021"
    }
}

# Every slide layout
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $lo = $p.SlideMaster.CustomLayouts.Item($i)
    $loShapes = $lo.Shapes
    for ($j = 1; $j -le $loShapes.Count; $j++) {
        $sh = $loShapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "11/1/2021"
        }
    }
}
